$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.220.30"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "2.175.13"
$ws.Range("E3").Value = "  -1.78%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.25"
$ws.Range("E5").Value = "  +5.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.87"

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.569"
$ws.Range("E9").Value = "  +2.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.12"
$ws.Range("E10").Value = "  -4.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.89"
$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0931"
$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.03"
$ws.Range("E13").Value = "  +4.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.104"
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").Value = "2.509.66"
$ws.Range("E15").Value = "  -1.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.868"
$ws.Range("E16").Value = "  +3.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.36"
$ws.Range("E17").Value = "  -2.91%  "

$ws.Range("D18").Value = "2.163.59"
$ws.Range("E18").Value = "  -2.15%  "

$ws.Range("D19").Value = "41.194.22"
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").Value = "0.0₃0953"
$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.13"
$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.72"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.11"
$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.01"
$ws.Range("E24").Value = "  -0.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.93"
$ws.Range("E25").Value = "  +9.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.68"
$ws.Range("E26").Value = "  +19.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.52"
$ws.Range("E28").Value = "  +4.35%  "

$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.54"
$ws.Range("E30").Value = "  -2.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.54"
$ws.Range("E31").Value = "  +0.62%  "

$ws.Range("E32").Value = "  -2.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0751"
$ws.Range("E33").Value = "  +5.67%  "

$ws.Range("E34").Value = "  -0.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.47"
$ws.Range("E35").Value = "  +4.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.50"
$ws.Range("E36").Value = "  +11.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.17"
$ws.Range("E37").Value = "  +6.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.58"
$ws.Range("E38").Value = "  -0.27%  "

$ws.Range("E39").Value = "  +6.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.18"
$ws.Range("E40").Value = "  -3.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.28"
$ws.Range("E41").Value = "  +15.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.65"
$ws.Range("E42").Value = "  -3.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.30"
$ws.Range("E43").Value = "  -1.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.04"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.198"
$ws.Range("E45").Value = "  +0.19%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.101"
$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.60"
$ws.Range("E47").Value = "  -0.67%  "

$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.18"
$ws.Range("E49").Value = "  +7.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.17"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.22"
$ws.Range("E51").Value = "  -5.98%  "
